$d = $word.ActiveDocument

$replacements = @(
    @("2025-11-24 Monday", "2025-11-25 Tuesday"),
    @("170×9=", "754×5="),
    @("864×8=", "863×9="),
    @("336×2=", "693×9="),
    @("530×4=", "802×2="),
    @("667×3=", "545×2="),
    @("185×2=", "531×4="),
    @("436×8=", "478×5="),
    @("467×6=", "436×7="),
    @("145×9=", "319×9="),
    @("320×9=", "952×6="),
    @("177×7=", "430×7="),
    @("279×2=", "898×8="),
    @("312×8=", "643×5="),
    @("849×6=", "288×7="),
    @("556×2=", "343×9="),
    @("699×5=", "408×5="),
    @("384×5=", "254×7="),
    @("987×5=", "108×7="),
    @("759×2=", "640×7="),
    @("599×5=", "675×5="),
    @("783×8=", "823×5="),
    @("421×5=", "549×6="),
    @("503×4=", "256×7="),
    @("712×5=", "236×2="),
    @("149×6=", "597×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
